$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the header style (bold, centered, bordered) from an existing header cell
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats

# New data cells
$ws.Range("G2").Value = 0.1239050709499376
$ws.Range("H2").Value = 0.991
